$d = $word.ActiveDocument

# Merge runs that were previously split apart only due to now-removed
# spelling/grammar proofing marks (w:proofErr), by re-applying a same-text
# Find & Replace across each affected span so Word collapses the runs.
$d.Content.Find.Execute("f some often used terms", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "f some often used terms", 2) | Out-Null
$d.Content.Find.Execute("devices for the control or support of conception;", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "devices for the control or support of conception;", 2) | Out-Null
$d.Content.Find.Execute("products specifically intended for the cleaning, disinfection or sterilisation ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "products specifically intended for the cleaning, disinfection or sterilisation ", 2) | Out-Null
$d.Content.Find.Execute("authorised by national law by virtue of that person's professional qualifications which gives, under that person's responsibility, specific design characteristics, and is intended for the sole use of a particular patient exclusively to", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "authorised by national law by virtue of that person's professional qualifications which gives, under that person's responsibility, specific design characteristics, and is intended for the sole use of a particular patient exclusively to", 2) | Out-Null
$d.Content.Find.Execute("written prescriptions of any authorised person shall not be considered to be custom-made devices", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "written prescriptions of any authorised person shall not be considered to be custom-made devices", 2) | Out-Null
$d.Content.Find.Execute("linical evaluation’ means a systematic and planned process to continuously generate, collect, analyse and assess", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "linical evaluation’ means a systematic and planned process to continuously generate, collect, analyse and assess", 2) | Out-Null
$d.Content.Find.Execute("monitoring, statistical considerations, organisation and conduct of a clinical investigation", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "monitoring, statistical considerations, organisation and conduct of a clinical investigation", 2) | Out-Null
$d.Content.Find.Execute("clinically relevant information coming from post-market surveillance, in particular the post-market clinical", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "clinically relevant information coming from post-market surveillance, in particular the post-market clinical", 2) | Out-Null
$d.Content.Find.Execute("relevant to the subject's decision to participate or, in the case of minors and of incapacitated subjects, an authorisation or agreement from their legally designated representative to include them in the clinical investigation", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "relevant to the subject's decision to participate or, in the case of minors and of incapacitated subjects, an authorisation or agreement from their legally designated representative to include them in the clinical investigation", 2) | Out-Null
$d.Content.Find.Execute("hospitalisation or prolongation of patient hospitalisation,", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "hospitalisation or prolongation of patient hospitalisation,", 2) | Out-Null
$d.Content.Find.Execute("foetal distress, foetal death or a congenital physical or mental impairment or birth defect", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "foetal distress, foetal death or a congenital physical or mental impairment or birth defect", 2) | Out-Null
$d.Content.Find.Execute("Human factors engineering and usability engineering can be considered to be", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Human factors engineering and usability engineering can be considered to be", 2) | Out-Null
$d.Content.Find.Execute("Regulation (EU) 2017/745 of the European Parliament and of the Council of 5 April 2017 on medical devices, amending Directive 2001/83/EC, Regulation (EC) No 178/2002 and Regulation (EC) No 1223/2009 and repealing Council Directives 90/385/EEC and 93/42/EEC", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Regulation (EU) 2017/745 of the European Parliament and of the Council of 5 April 2017 on medical devices, amending Directive 2001/83/EC, Regulation (EC) No 178/2002 and Regulation (EC) No 1223/2009 and repealing Council Directives 90/385/EEC and 93/42/EEC", 2) | Out-Null
$d.Content.Find.Execute("IEC 62366-1 Medical devices – Part 1: Application of usability engineering to medical devices. Edition 1.0. ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "IEC 62366-1 Medical devices – Part 1: Application of usability engineering to medical devices. Edition 1.0. ", 2) | Out-Null
$d.Content.Find.Execute("Geneva: International Electrotectnical Commission; 2015. ISBN:978-2-8322-2281-2.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Geneva: International Electrotectnical Commission; 2015. ISBN:978-2-8322-2281-2.", 2) | Out-Null
$d.Content.Find.Execute("FDA (2016) Applying Human Factors and Usability Engineering to Medical Devices: Guidance for Industry and Food and Drug Administration Staff.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "FDA (2016) Applying Human Factors and Usability Engineering to Medical Devices: Guidance for Industry and Food and Drug Administration Staff.", 2) | Out-Null

# Remove obsolete bracketed reference markers that were left over from
# footnote-style citations no longer used in the document.
$d.Content.Find.Execute("EU Regulation 2017/745 [1]", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "EU Regulation 2017/745", 2) | Out-Null
$d.Content.Find.Execute("IEC 62366-1 [2]", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "IEC 62366-1", 2) | Out-Null
$d.Content.Find.Execute("FDA Guidance [3]", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "FDA Guidance", 2) | Out-Null
